# Auto-generated edit script: update recalculated market-data-driven
# profit columns (H-N) across several worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2382.1155
$ws.Range("J131").Value = 3272.4614
$ws.Range("L131").Value = 9817.3842
$ws.Range("N131").Value = -19897.3842
$ws.Range("H132").Value = 15162.939
$ws.Range("I132").Value = 1993.6167
$ws.Range("J132").Value = 146856.17
$ws.Range("K132").Value = 5980.8501
$ws.Range("L132").Value = 440568.51
$ws.Range("M132").Value = -3450.8501
$ws.Range("N132").Value = -445628.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3329.0908
$ws.Range("I63").Value = 2034.2858
$ws.Range("J63").Value = 5595
$ws.Range("K63").Value = 2034.2858
$ws.Range("L63").Value = 5595
$ws.Range("M63").Value = -1348.2858
$ws.Range("N63").Value = -6967
$ws.Range("H66").Value = 3329.0908
$ws.Range("I66").Value = 2034.2858
$ws.Range("J66").Value = 5595
$ws.Range("K66").Value = 10171.429
$ws.Range("L66").Value = 27975
$ws.Range("M66").Value = -6739.429
$ws.Range("N66").Value = -34839
$ws.Range("H123").Value = 38999.5
$ws.Range("J123").Value = 38999.5
$ws.Range("L123").Value = 38999.5
$ws.Range("N123").Value = -48799.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 9724.75
$ws.Range("J81").Value = 9724.75
$ws.Range("L81").Value = 9724.75
$ws.Range("N81").Value = -11846.75
$ws.Range("H84").Value = 9724.75
$ws.Range("J84").Value = 9724.75
$ws.Range("L84").Value = 29174.25
$ws.Range("N84").Value = -39782.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4390327
$ws.Range("I31").Value = 1802.0646
$ws.Range("K31").Value = 1802.0646
$ws.Range("M31").Value = -1507.0646
$ws.Range("H34").Value = 4390327
$ws.Range("I34").Value = 1802.0646
$ws.Range("K34").Value = 1802.0646
$ws.Range("M34").Value = -1600.0646
$ws.Range("H58").Value = 1491.4108
$ws.Range("I58").Value = 1086.7428
$ws.Range("J58").Value = 2165.8572
$ws.Range("K58").Value = 1086.7428
$ws.Range("L58").Value = 2165.8572
$ws.Range("M58").Value = -883.7428
$ws.Range("N58").Value = -2571.8572
$ws.Range("H132").Value = 63342.566
$ws.Range("I132").Value = 1725.6
$ws.Range("J132").Value = 178874.38
$ws.Range("K132").Value = 5176.799999999999
$ws.Range("L132").Value = 536623.14
$ws.Range("M132").Value = -2646.799999999999
$ws.Range("N132").Value = -541683.14
$ws.Range("H136").Value = 1491.4108
$ws.Range("I136").Value = 1086.7428
$ws.Range("J136").Value = 2165.8572
$ws.Range("K136").Value = 3260.2284
$ws.Range("L136").Value = 6497.571599999999
$ws.Range("M136").Value = -710.2284
$ws.Range("N136").Value = -11597.5716
$ws.Range("H137").Value = 27711.111
$ws.Range("J137").Value = 27711.111
$ws.Range("L137").Value = 27711.111
$ws.Range("N137").Value = -37911.111
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5251.375
$ws.Range("J70").Value = 5999.8335
$ws.Range("L70").Value = 17999.5005
$ws.Range("N70").Value = -18629.5005
$ws.Range("H73").Value = 5251.375
$ws.Range("J73").Value = 5999.8335
$ws.Range("L73").Value = 17999.5005
$ws.Range("N73").Value = -20183.5005
$ws.Range("H74").Value = 13869
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13869
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 41607
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -43729
$ws.Range("H75").Value = 2973.9092
$ws.Range("I75").Value = 606.5
$ws.Range("J75").Value = 3500
$ws.Range("K75").Value = 1819.5
$ws.Range("L75").Value = 10500
$ws.Range("M75").Value = -821.5
$ws.Range("N75").Value = -12496
$ws.Range("H76").Value = 3056.5
$ws.Range("I76").Value = 1213
$ws.Range("J76").Value = 4900
$ws.Range("K76").Value = 3639
$ws.Range("L76").Value = 14700
$ws.Range("M76").Value = -3256
$ws.Range("N76").Value = -15466
$ws.Range("H77").Value = 13869
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13869
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 124821
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -135429
$ws.Range("H78").Value = 2973.9092
$ws.Range("I78").Value = 606.5
$ws.Range("J78").Value = 3500
$ws.Range("K78").Value = 5458.5
$ws.Range("L78").Value = 31500
$ws.Range("M78").Value = -466.5
$ws.Range("N78").Value = -41484
$ws.Range("H79").Value = 3056.5
$ws.Range("I79").Value = 1213
$ws.Range("J79").Value = 4900
$ws.Range("K79").Value = 3639
$ws.Range("L79").Value = 14700
$ws.Range("M79").Value = -2313
$ws.Range("N79").Value = -17352
$ws.Range("H80").Value = 91116170
$ws.Range("I80").Value = 2000000
$ws.Range("J80").Value = 100027784
$ws.Range("K80").Value = 6000000
$ws.Range("L80").Value = 300083352
$ws.Range("M80").Value = -5999064
$ws.Range("N80").Value = -300085224
$ws.Range("H81").Value = 1763
$ws.Range("I81").Value = 933.3333
$ws.Range("J81").Value = 3007.5
$ws.Range("K81").Value = 2799.9999
$ws.Range("L81").Value = 9022.5
$ws.Range("M81").Value = -1676.9999
$ws.Range("N81").Value = -11268.5
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 91116170
$ws.Range("I83").Value = 2000000
$ws.Range("J83").Value = 100027784
$ws.Range("K83").Value = 18000000
$ws.Range("L83").Value = 900250056
$ws.Range("M83").Value = -17995320
$ws.Range("N83").Value = -900259416
$ws.Range("H84").Value = 1763
$ws.Range("I84").Value = 933.3333
$ws.Range("J84").Value = 3007.5
$ws.Range("K84").Value = 8399.9997
$ws.Range("L84").Value = 27067.5
$ws.Range("M84").Value = -2783.9997
$ws.Range("N84").Value = -38299.5
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H87").Value = 24285.715
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H88").Value = 5977.6875
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5977.6875
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 17933.0625
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -18789.0625
$ws.Range("H90").Value = 24285.715
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H91").Value = 5977.6875
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5977.6875
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 17933.0625
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -20897.0625
$ws.Range("H121").Value = 89420.23
$ws.Range("I121").Value = 310
$ws.Range("J121").Value = 111697.79
$ws.Range("K121").Value = 930
$ws.Range("L121").Value = 335093.37
$ws.Range("M121").Value = 380
$ws.Range("N121").Value = -337713.37
$ws.Range("H131").Value = 867.62244
$ws.Range("J131").Value = 872.36456
$ws.Range("L131").Value = 2617.09368
$ws.Range("N131").Value = -12697.09368

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 42727.6
$ws.Range("J137").Value = 42727.6
$ws.Range("L137").Value = 42727.6
$ws.Range("N137").Value = -52927.6
$ws.Range("H139").Value = 56742.145
$ws.Range("I139").Value = 160000
$ws.Range("J139").Value = 39532.5
$ws.Range("K139").Value = 160000
$ws.Range("L139").Value = 39532.5
$ws.Range("M139").Value = -154860
$ws.Range("N139").Value = -49812.5
$ws.Range("H141").Value = 38013.57
$ws.Range("J141").Value = 38013.57
$ws.Range("L141").Value = 38013.57
$ws.Range("N141").Value = -48373.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1216.9286
$ws.Range("I122").Value = 1214.6666
$ws.Range("J122").Value = 1221
$ws.Range("K122").Value = 3643.9998
$ws.Range("L122").Value = 3663
$ws.Range("M122").Value = -1193.9998
$ws.Range("N122").Value = -8563
$ws.Range("H135").Value = 40357.145
$ws.Range("J135").Value = 40357.145
$ws.Range("L135").Value = 40357.145
$ws.Range("N135").Value = -50497.145

